# Fix lỗi update data
# Adds a new "Sheet2" worksheet containing employee login data
# (Tài khoản / Mật Khẩu / Họ Tên / CCCD) after Sheet1, and restores
# the previously-selected cell on Sheet1.

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Sheet1")

# --- add the new worksheet, placed right after Sheet1 -------------------
$ws2 = $wb.Worksheets.Add($null, $ws1)
$ws2.Name = "Sheet2"

# CCCD column (D) must stay text, same as on Sheet1, so the leading
# zeros in the id numbers aren't lost.
$ws2.Range("D2:D5").NumberFormat = "@"

# --- header row (column order matches the data-entry order used when
#     the sheet was authored: header first, then row 2 in full, then
#     columns A/B for rows 3-5, then column C for rows 3-5) -------------
$ws2.Range("A1").Value = "Tài khoản"
$ws2.Range("B1").Value = "Mật Khẩu"
$ws2.Range("C1").Value = "Họ Tên"
$ws2.Range("D1").Value = "CCCD"

$ws2.Range("A2").Value = "nv1"
$ws2.Range("B2").Value = "nv1"
$ws2.Range("C2").Value = "Nhân viên 1"
$ws2.Range("D2").Value = "000000000001"

$ws2.Range("A3").Value = "nv2"
$ws2.Range("B3").Value = "nv2"
$ws2.Range("A4").Value = "nv3"
$ws2.Range("B4").Value = "nv3"
$ws2.Range("A5").Value = "nv4"
$ws2.Range("B5").Value = "nv4"

$ws2.Range("C3").Value = "Nhân viên 2"
$ws2.Range("C4").Value = "Nhân viên 3"
$ws2.Range("C5").Value = "Nhân viên 4"

$ws2.Range("D3").Value = "000000000002"
$ws2.Range("D4").Value = "000000000003"
$ws2.Range("D5").Value = "000000000004"

# --- header formatting: bold, centered -----------------------------------
# (format A1 fully, then copy its format across the rest of the header so
#  only the final bold+centered style ends up in the saved style table)
$a1 = $ws2.Range("A1")
$a1.Font.Bold = $true
$a1.HorizontalAlignment = -4108  # xlCenter
$a1.Copy()
$ws2.Range("B1:D1").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

# --- column widths (auto-fit to content, like the original bestFit) -----
$ws2.Columns.Item("B:D").AutoFit()

# --- view settings for the new sheet ------------------------------------
$ws2.Range("A2:D5").Select()
$excel.ActiveWindow.Zoom = 220

# --- Sheet1's view: scrolled right a couple columns, G3 active ----------
$ws1.Activate()
$ws1.Application.ActiveWindow.ScrollColumn = 3
$ws1.Range("G3").Select()
